# Update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "57.882.11"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "3.148.10"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  +1.20%  "

$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "530.94"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "140.54"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "

$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "3.149.26"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  +1.17%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "0.445"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  +2.87%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "7.20"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  -0.95%  "

$ws.Range("E11").Value = "  -0.56%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "0.399"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  +3.83%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "3.694.32"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  +1.32%  "

$ws.Range("E14").Value = "  +3.04%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "25.67"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -2.08%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0000165"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "58.038.93"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "3.151.26"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  +1.59%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "6.13"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "12.90"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  +0.54%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "7.98"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "357.05"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +5.99%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "0.999"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("E24").Value = "  +3.47%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "0.511"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("E26").Value = "  +1.41%  "

$ws.Range("E27").Value = "  +0.29%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0933"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  +0.52%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "7.46"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  +3.13%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "6.38"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -3.26%  "

$ws.Range("E32").Value = "  +1.47%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "21.24"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +1.51%  "

$ws.Range("E34").Value = "  -1.04%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "4.94"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +5.83%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "158.04"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +2.47%  "

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "6.20"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "26.13"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -3.34%  "

$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "1.28"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  -1.65%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0672"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  +0.61%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "1.63"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +9.33%  "

$ws.Range("E42").Value = "  +5.86%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "0.706"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  +3.03%  "

$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "3.193.02"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  +1.30%  "

$ws.Range("E45").Value = "  +6.01%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "36.77"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -0.50%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "2.335.09"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  +1.62%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "0.993"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "6.08"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +1.40%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "20.37"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -2.57%  "

